$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: C2, D2, F2 change (E2 stays 2)
$ws.Cells.Item(2, 3).Value = "20"
$ws.Cells.Item(2, 4).Value = "12"
$ws.Cells.Item(2, 6).Value = "1"

# Row 3: C3, D3, E3 change (F3 stays 0)
$ws.Cells.Item(3, 3).Value = "6"
$ws.Cells.Item(3, 4).Value = "4"
$ws.Cells.Item(3, 5).Value = "1"

# Row 4: C4, D4, E4 change (F4 stays 0)
$ws.Cells.Item(4, 3).Value = "0"
$ws.Cells.Item(4, 4).Value = "7"
$ws.Cells.Item(4, 5).Value = "0"

# Row 5: C5, D5, E5, F5 all change
$ws.Cells.Item(5, 3).Value = "1"
$ws.Cells.Item(5, 4).Value = "2"
$ws.Cells.Item(5, 5).Value = "0"
$ws.Cells.Item(5, 6).Value = "0"

# Row 6: C6, D6, E6, F6 all change
$ws.Cells.Item(6, 3).Value = "42"
$ws.Cells.Item(6, 4).Value = "26"
$ws.Cells.Item(6, 5).Value = "2"
$ws.Cells.Item(6, 6).Value = "2"

# Row 7: C7, D7, E7 change (F7 stays 0)
$ws.Cells.Item(7, 3).Value = "16"
$ws.Cells.Item(7, 4).Value = "18"
$ws.Cells.Item(7, 5).Value = "1"

# Row 8: only D8 changes
$ws.Cells.Item(8, 4).Value = "2"

# Row 9: C9, D9, E9 change (F9 stays 0)
$ws.Cells.Item(9, 3).Value = "1"
$ws.Cells.Item(9, 4).Value = "6"
$ws.Cells.Item(9, 5).Value = "0"
